$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# Sheet "Summary": update running totals after Trade #7 closed
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")
$summary.Range("B3").Value = 1199.85   # Current Capital
$summary.Range("B4").Value = -0.14     # Total P&L $
$summary.Range("B5").Value = -0.4      # Total P&L %
$summary.Range("B6").Value = 7         # Total Trades
$summary.Range("B8").Value = 6         # Losing Trades
$summary.Range("B9").Value = 14.29     # Win Rate %

# ---------------------------------------------------------------------------
# Sheet "Strategy Status": update the MarketMaking strategy row (row 4)
# ---------------------------------------------------------------------------
$status = $wb.Worksheets.Item("Strategy Status")
$status.Range("C4").Value = 99.84999999999999  # Capital
$status.Range("D4").Value = 7                  # Trades
$status.Range("E4").Value = -0.14              # P&L $
$status.Range("F4").Value = -0.15              # P&L %
$status.Range("G4").Value = 14.29              # Win Rate %

# ---------------------------------------------------------------------------
# Sheets "All Trades" and "MarketMaking": append the newly closed trade #7
# as row 8 on both sheets (they mirror each other).
# ---------------------------------------------------------------------------
function Add-TradeRow8($ws) {
    $ws.Range("A8").Value = 7
    # B8/C8 look like a date / a time; force them to remain plain text so
    # they aren't auto-converted to date/time serial numbers.
    $ws.Range("B8").NumberFormat = "@"
    $ws.Range("B8").Value = "2026-02-17"
    $ws.Range("B8").Style = "Normal"
    $ws.Range("C8").NumberFormat = "@"
    $ws.Range("C8").Value = "13:33:48"
    $ws.Range("C8").Style = "Normal"
    $ws.Range("D8").Value = "MarketMaking"
    $ws.Range("E8").Value = "UP"
    $ws.Range("F8").Value = 0.07000000000000001
    $ws.Range("G8").Value = 0.028801
    $ws.Range("H8").Value = "CLOSED"
    $ws.Range("I8").Value = -58.8561
    $ws.Range("J8").Value = -0.04
    $ws.Range("K8").Value = 99.84999999999999
    $ws.Range("L8").Value = 0
    $ws.Range("M8").Value = 0
    $ws.Range("N8").Value = 0.6
    $ws.Range("O8").Value = "Normal spread capture: 19600 bps"
    $ws.Range("P8").Value = "early_exit"
    $ws.Range("Q8").Value = 0.13
}

$allTrades = $wb.Worksheets.Item("All Trades")
Add-TradeRow8 $allTrades

$marketMaking = $wb.Worksheets.Item("MarketMaking")
Add-TradeRow8 $marketMaking
